# Generate Report for Handoff
#
# Two new files have been handed off since the last report:
#   - 26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png  (depends on the .md file below)
#   - cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png   (depends on the .md file below)
#   - da621f65-d76e-4e0e-ac04-3d89b169c3c4.md    (the dependency itself)
# replacing the old 3f7af0e6-3842-406e-91a3-430a4d9c9fb8.md row.
#
# This script appends the two new rows to the "Overview" sheet and to each
# per-locale sheet ("zh-cn", "de-de"), refreshes the existing first-row
# values to point at the new handed-off file, and re-creates the hyperlinks.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # matches the workbook's existing HyperLink font (RGB 0x64,0x95,0xED)

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-51-11 10:51:13"

$ov.Range("A3").Value = "cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-51-11 10:51:13"

$ov.Range("A4").Value = "da621f65-d76e-4e0e-ac04-3d89b169c3c4.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-51-11 10:51:13"

Style-AsHyperlink($ov.Range("A2"))
Style-AsHyperlink($ov.Range("A3"))
Style-AsHyperlink($ov.Range("A4"))

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png", "", "", "26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png", "", "", "cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/da621f65-d76e-4e0e-ac04-3d89b169c3c4.md", "", "", "da621f65-d76e-4e0e-ac04-3d89b169c3c4.md") | Out-Null

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn" and "de-de")
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Row2Time = "2016-03-11 10:51:10"; Row3Time = "2016-03-11 10:51:10"; Row4Time = "2016-03-11 10:51:10"; XlfD2 = "99fba8e063f50c148e1f3dba5d789d82eb1ef9ee.png"; XlfD3 = "b3aea06de7ba5f2185e73ea0a875f1de1b610de7.png"; XlfD4 = "da621f65-d76e-4e0e-ac04-3d89b169c3c4.19cb4ae79004bcc204757e1938bb92b12eb930b0.zh-cn.xlf" },
    @{ Name = "de-de"; Row2Time = "2016-03-11 10:51:13"; Row3Time = "2016-03-11 10:51:13"; Row4Time = "2016-03-11 10:51:13"; XlfD2 = "99fba8e063f50c148e1f3dba5d789d82eb1ef9ee.png"; XlfD3 = "b3aea06de7ba5f2185e73ea0a875f1de1b610de7.png"; XlfD4 = "da621f65-d76e-4e0e-ac04-3d89b169c3c4.19cb4ae79004bcc204757e1938bb92b12eb930b0.de-de.xlf" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)

    # Row 2 - first .png (is itself a dependency of the .md file)
    $ws.Range("A2").Value = "26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png"
    $ws.Range("B2").Value = ".png"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("D2").Value = $loc.XlfD2
    $ws.Range("E2").Value = $loc.Row2Time
    $ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "IsDependency"
    $ws.Range("J2").Value = "e2e\da621f65-d76e-4e0e-ac04-3d89b169c3c4.md"

    # Row 3 - second .png (also a dependency of the .md file)
    $ws.Range("A3").Value = "cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png"
    $ws.Range("B3").Value = ".png"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $loc.XlfD3
    $ws.Range("E3").Value = $loc.Row3Time
    $ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "IsDependency"
    $ws.Range("J3").Value = "e2e\da621f65-d76e-4e0e-ac04-3d89b169c3c4.md"

    # Row 4 - the .md file itself (included directly, not a dependency)
    $ws.Range("A4").Value = "da621f65-d76e-4e0e-ac04-3d89b169c3c4.md"
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Range("D4").Value = $loc.XlfD4
    $ws.Range("E4").Value = $loc.Row4Time
    $ws.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "Include"

    Style-AsHyperlink($ws.Range("A2"))
    Style-AsHyperlink($ws.Range("B2"))
    Style-AsHyperlink($ws.Range("D2"))
    Style-AsHyperlink($ws.Range("A3"))
    Style-AsHyperlink($ws.Range("B3"))
    Style-AsHyperlink($ws.Range("D3"))
    Style-AsHyperlink($ws.Range("A4"))
    Style-AsHyperlink($ws.Range("B4"))
    Style-AsHyperlink($ws.Range("D4"))

    $base = "https://github.com/OpenLocalizationTest/oltest/blob/e2e/"
    $ws.Hyperlinks.Add($ws.Range("A2"), ($base + "26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png"), "", "", "26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), ($base + "26831a4e-66d1-4b83-9b3c-c53fa18eb65e.png"), "", "", ".png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/" + $loc.XlfD2), "", "", $loc.XlfD2) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), ($base + "cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png"), "", "", "cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), ($base + "cc2027e7-d036-4f52-80bc-d47e2f68ed3a.png"), "", "", ".png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/" + $loc.XlfD3), "", "", $loc.XlfD3) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), ($base + "da621f65-d76e-4e0e-ac04-3d89b169c3c4.md"), "", "", "da621f65-d76e-4e0e-ac04-3d89b169c3c4.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), ($base + "da621f65-d76e-4e0e-ac04-3d89b169c3c4.md"), "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D4"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/" + $loc.XlfD4), "", "", $loc.XlfD4) | Out-Null
}

Write-Host "Report regenerated for handoff."
